$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "they can be used in entertainment industry, for example for making movies or games.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "they can be used in the entertainment industry, for example in making movies or games.",
    2
)
